$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.996.40"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.81%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.879.31"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.95%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9979"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.38"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -4.50%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9980"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.13%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.79%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2926"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.91%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06625"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.07%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.877.83"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.01%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "16.74"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -4.01%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07248"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6677"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.78%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "86.33"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.81%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.877"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.82%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.953.90"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.95%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007921"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.41%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9980"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.17%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.78"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.120.69"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.92%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9974"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.13%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.759"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.17%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.648"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.61%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.059"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.40%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.15"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.79%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.26"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.41%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.14"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.92%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.910"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -5.65%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.388"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.37%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.174"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.70%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08779"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.79%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.945"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.18%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05068"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.57%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7106"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.05%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.109"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.74%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.665"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.77%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.696"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -5.08%  "

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.68%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.182"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -5.81%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9303"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.87%  "

# Row 41
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4253"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.66%  "

# Row 42
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.779"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -5.45%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9976"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.14%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.00"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.08%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.461"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -3.50%  "

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.16%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05650"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -1.73%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "32.42"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.20%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.3756"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.23%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.221"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.59%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "55.82"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.02%  "
